# ---------------------------------------------------------------------------
# Adds a new "2022-Q3" quarter sheet (copied/derived from the "2022-Q2"
# layout) right after "总计", renumbers nothing else (Excel COM shifts the
# remaining sheets automatically), and updates the "总计" summary sheet with
# the new quarter's totals (existing rows shift down by one row).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: write a text value into a cell while defeating Excel's "this looks
# like a number" auto-conversion (needed for fund codes like "000689" and for
# the numeric-looking percentage/NAV strings that must stay text), then
# restore the cell's original (unstyled) appearance so formatting is not
# disturbed by the NumberFormat round-trip.
function Set-TextValue($cell, [string]$value) {
    $savedStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $savedStyle
}

# ===========================================================================
# 1. Create the new "2022-Q3" sheet by copying the "2022-Q2" sheet's layout
#    (keeps headers / fonts / borders identical) right after "总计".
# ===========================================================================
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($null, $totalSheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# The copied sheet has 17 data rows (rows 2-18); the new quarter only has 9
# data rows (rows 2-10), so remove the extra rows entirely.
$q3Sheet.Range("A11:H18").Clear()

# ---------------------------------------------------------------------------
# Fill in the "2022-Q3" data (row 1 header + row-index column A already
# carry over correctly from the copy).
# ---------------------------------------------------------------------------
$q3Data = @(
    @("000689", "前海开源新经济灵活配置混合A", "98.71", "92.74", "4.39", "4.3334", 9),
    @("010490", "鹏华高质量增长混合A",        "12.74", "93.98", "9.32", "1.1874", 2),
    @("013157", "前海开源新经济灵活配置混合C", "14.05", "92.74", "4.39", "0.6168", 9),
    @("009023", "鹏华稳健回报混合",            "4.12",  "94.39", "9.60", "0.3955", 1),
    @("005314", "万家中证1000指数增强C",       "14.28", "94.11", "1.01", "0.1442", 5),
    @("005313", "万家中证1000指数增强A",       "13.25", "94.11", "1.01", "0.1338", 5),
    @("010491", "鹏华高质量增长混合C",        "0.44",  "93.98", "9.32", "0.0410", 2),
    @("013489", "广发东财大数据精选灵活配置混合C", "0.18", "60.68", "0.88", "0.0016", 8),
    @("002802", "广发东财大数据精选灵活配置混合A", "0.15", "60.68", "0.88", "0.0013", 8)
)

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $row = $i + 2
    $rec = $q3Data[$i]

    Set-TextValue $q3Sheet.Cells.Item($row, 2) $rec[0]   # B: 基金代码
    $q3Sheet.Cells.Item($row, 3).Value = $rec[1]          # C: 基金名称 (plain text already)
    Set-TextValue $q3Sheet.Cells.Item($row, 4) $rec[2]   # D: 基金规模
    Set-TextValue $q3Sheet.Cells.Item($row, 5) $rec[3]   # E: 股票总仓位
    Set-TextValue $q3Sheet.Cells.Item($row, 6) $rec[4]   # F: 仓位占比
    Set-TextValue $q3Sheet.Cells.Item($row, 7) $rec[5]   # G: 持有市值(亿元)
    $q3Sheet.Cells.Item($row, 8).Value = $rec[6]          # H: 仓位排名 (real number)
}

# ===========================================================================
# 2. Update the "总计" (totals) sheet: the existing rows (2022-Q2 .. 2020-Q4)
#    shift down by one row, and a brand-new row 2 for "2022-Q3" is added.
#    Column A (row index) is left untouched on the existing rows and simply
#    extended with the next sequential value for the appended row.
# ===========================================================================

# Shift existing rows 8->9, 7->8, ..., 2->3 (bottom-up so we never overwrite
# data before it has been read).
for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    $totalSheet.Cells.Item($r, 1).Copy($totalSheet.Cells.Item($dest, 1))
    $totalSheet.Cells.Item($dest, 1).Value = $dest - 2

    $totalSheet.Cells.Item($dest, 2).Value = $totalSheet.Cells.Item($r, 2).Value
    $totalSheet.Cells.Item($dest, 3).Value = $totalSheet.Cells.Item($r, 3).Value
    $totalSheet.Cells.Item($dest, 4).Value = $totalSheet.Cells.Item($r, 4).Value
}

# Now write the brand-new "2022-Q3" row into row 2.
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 9
$totalSheet.Cells.Item(2, 4).Value = 6.86

Write-Output "2022-Q3 sheet inserted and 总计 sheet updated."
